# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.235.67"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "2.226.71"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'244.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").Value = "'0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("D7").Value = "'73.87"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.88%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").Value = "'0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("D10").Value = "'42.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.73%  "

$ws.Range("D11").Value = "'0.0974"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.12%  "

$ws.Range("D12").Value = "'7.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.80%  "

$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").Value = "'14.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("D15").Value = "'0.851"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("D16").Value = "2.232.90"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("E17").Value = "  +19.62%  "

$ws.Range("D18").Value = "42.108.27"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").Value = "'6.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.19%  "

$ws.Range("D20").Value = "'72.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.65%  "

$ws.Range("D21").Value = "'10.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +38.73%  "

$ws.Range("D22").Value = "'230.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("E23").Value = "  -3.65%  "

$ws.Range("D24").Value = "'11.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.28%  "

$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("E26").Value = "  -1.23%  "

$ws.Range("D27").Value = "'2.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "

$ws.Range("E28").Value = "  +3.22%  "

$ws.Range("D29").Value = "'167.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("E30").Value = "  +3.56%  "

$ws.Range("D31").Value = "'5.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +18.58%  "

$ws.Range("D32").Value = "'0.0802"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").Value = "'29.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").Value = "'4.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.78%  "

$ws.Range("E37").Value = "  +3.25%  "

$ws.Range("E38").Value = "  -0.69%  "

$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("D40").Value = "'5.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("D41").Value = "'62.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.15%  "

$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("D43").Value = "'8.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("D44").Value = "'105.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.74%  "

$ws.Range("E45").Value = "  +3.25%  "

$ws.Range("E47").Value = "  +7.80%  "

$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").Value = "'1.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.62%  "

$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").Value = "'4.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
